$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D (Price) values: force text format to preserve exact
# string representation (e.g. trailing zeros, multi-dot formatted
# numbers) since these are display strings, not real numbers. ---
$priceCells = @{
    'D2' = '29.464.29'
    'D3' = '1.853.23'
    'D4' = '0.9991'
    'D6' = '0.6305'
    'D7' = '1.000'
    'D8' = '0.07613'
    'D9' = '0.2978'
    'D10' = '24.56'
    'D11' = '0.07712'
    'D12' = '1.921.74'
    'D13' = '5.007'
    'D14' = '0.6886'
    'D15' = '83.47'
    'D16' = '0.000009919'
    'D17' = '2.188.45'
    'D18' = '6.174'
    'D19' = '29.587.12'
    'D20' = '232.98'
    'D21' = '12.54'
    'D22' = '7.706'
    'D24' = '1.001'
    'D25' = '154.85'
    'D26' = '0.1400'
    'D27' = '8.478'
    'D28' = '17.70'
    'D29' = '1.474'
    'D30' = '0.05799'
    'D32' = '4.128'
    'D33' = '4.020'
    'D36' = '0.7212'
    'D38' = '1.252.51'
    'D39' = '2.796'
    'D40' = '0.01802'
    'D41' = '0.9092'
    'D42' = '6.110'
    'D43' = '2.092.21'
    'D44' = '0.9998'
    'D45' = '67.81'
    'D46' = '101.69'
    'D47' = '7.325'
    'D48' = '0.00000000118'
    'D49' = '9.163'
    'D50' = '0.4040'
    'D51' = '1.710'
}
foreach ($addr in $priceCells.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $priceCells[$addr]
}

# --- Coin name / link / volume columns (plain text, no numeric
# auto-conversion risk) ---
$textCells = @{
    'E2' = '  -1.01%  '
    'E3' = '  -0.07%  '
    'E4' = '  -0.11%  '
    'E6' = '  -3.33%  '
    'E7' = '  -0.06%  '
    'E8' = '  +1.29%  '
    'E9' = '  -0.08%  '
    'E10' = '  -0.25%  '
    'E11' = '  +1.01%  '
    'E12' = '  +3.61%  '
    'B13' = 'Polkadot'
    'C13' = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
    'E13' = '  -0.97%  '
    'B14' = 'Polygon'
    'C14' = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
    'E14' = '  +0.24%  '
    'E15' = '  -0.18%  '
    'E16' = '  +1.97%  '
    'E17' = '  +3.74%  '
    'E18' = '  +1.13%  '
    'E19' = '  -0.73%  '
    'E20' = '  -1.48%  '
    'E21' = '  -0.95%  '
    'E22' = '  -0.95%  '
    'E23' = '  -0.03%  '
    'E24' = '  -0.09%  '
    'E25' = '  -2.30%  '
    'E26' = '  -1.98%  '
    'E27' = '  -0.96%  '
    'E28' = '  -0.94%  '
    'E29' = '  -1.19%  '
    'E30' = '  -4.53%  '
    'E31' = '  -1.64%  '
    'E32' = '  -0.52%  '
    'E33' = '  -1.42%  '
    'E34' = '  +1.07%  '
    'E36' = '  -0.27%  '
    'E37' = '  -0.73%  '
    'E38' = '  +4.14%  '
    'E39' = '  -0.47%  '
    'E40' = '  +0.58%  '
    'E41' = '  -0.32%  '
    'E42' = '  -2.99%  '
    'E43' = '  +3.18%  '
    'E44' = '  -0.10%  '
    'E45' = '  +1.38%  '
    'E46' = '  +0.30%  '
    'E47' = '  +0.91%  '
    'B48' = 'BabyDogeCoin'
    'C48' = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
    'E48' = '  -0.95%  '
    'B49' = 'EnergySwap'
    'C49' = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
    'E49' = '  -0.12%  '
    'B50' = 'TheSandbox'
    'C50' = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
    'E50' = '  -0.52%  '
    'E51' = '  +2.01%  '
}
foreach ($addr in $textCells.Keys) {
    $ws.Range($addr).Value = $textCells[$addr]
}
